# yolo映射表.xlsx — "Add files via upload"
#
# The dataset/id/name lookup table on Sheet1 gained two new rows:
#   - A new "NWPU_VHR-10" row (id 0 -> "airplane") inserted at the top of the
#     existing NWPU_VHR-10 block (old row 35), pushing that block (and
#     everything below it) down by one row.
#   - A new "VisDrone" row (id 0 -> "pedestrian") inserted at the top of the
#     existing VisDrone block (old row 44, row 45 after the first insert),
#     pushing that block (and everything below it) down by one row again.
#
# Net effect: dimension grows from A1:C73 to A1:C75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "NWPU_VHR-10" / 0 / "airplane" above the NWPU_VHR-10 block ---
$ws.Rows(35).Insert()
$ws.Range("A35").Value = "NWPU_VHR-10"
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = "airplane"
# Match the formatting used by the rest of the block / the row it was
# copied in front of (col A keeps the block's left-aligned style, col C
# picks up the "vertical-center" style carried down from row 34).
$ws.Range("A36").Copy()
$ws.Range("A35").PasteSpecial(-4122)

# --- Insert "VisDrone" / 0 / "pedestrian" above the VisDrone block ---
# (VisDrone's first row is now at 45, after the insert above.)
$ws.Rows(45).Insert()
$ws.Range("A45").Value = "VisDrone"
$ws.Range("B45").Value = 0
$ws.Range("C45").Value = "pedestrian"
# This block uses the plain default style throughout, matching the row
# it was inserted in front of.
$ws.Range("A46").Copy()
$ws.Range("A45").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Leave the selection where the editor ended up after typing the new data.
$ws.Range("C35").Select()
